$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# "Prueba 6" section is being repurposed from a "duplicate movie" scenario to
# a "delete movie" scenario. All edits below are scoped to that section via
# explicit paragraph indices (41-48 before any insertion) to avoid touching
# similarly-worded text earlier in the document (e.g. "Luego hace click...").
# ---------------------------------------------------------------------------

# Paragraph 42 (italic scenario summary under "Prueba 6")
# "El administrador no es capaz de editar la película que desea porque los
#  datos colisionan con los de una película ya creada."
#  -> "El administrador desea eliminar una película de la base de datos."
$p42 = $d.Paragraphs(42)
$r42 = $p42.Range
$r42.MoveEnd(1, -1) | Out-Null
$r42.Text = "El administrador desea eliminar una película de la base de datos."

# Paragraph 44 (Criterio de aceptación bullet)
# "La aplicación le debe preguntar al administrador si está seguro que desea
#  guardar los cambios de la película ya que ya existe una con esos datos."
#  -> "La película no debe aparecer más en la lista de la cartelera."
$p44 = $d.Paragraphs(44)
$r44 = $p44.Range
$r44.MoveEnd(1, -1) | Out-Null
$r44.Text = "La película no debe aparecer más en la lista de la cartelera."

# Paragraph 46 (first "Pasos" bullet)
# "El administrador llena los campos con los datos de la película. El nombre
#  de la película es idéntico al de una ya creada."
#  -> "El administrador entra al menú de cartelera."
$p46 = $d.Paragraphs(46)
$r46 = $p46.Range
$r46.MoveEnd(1, -1) | Out-Null
$r46.Text = "El administrador entra al menú de cartelera."

# Paragraph 47 (second "Pasos" bullet) keeps its middle "click" run (and the
# spell-check proofErr markers around it) untouched; only the text before and
# after "click" changes.
#   "Luego hace "  -> "Luego navega por las páginas de la misma hasta elegir
#                       una película y hace "
#   " en el botón "Guardar cambios"."
#                  -> " en el botón de eliminar que está al lado del botón de
#                       edición."
$p47 = $d.Paragraphs(47)
$pStart47 = $p47.Range.Start
$beforeClick = $d.Range($pStart47, $pStart47 + 11)
$beforeClick.Text = "Luego navega por las páginas de la misma hasta elegir una película y hace "

$full47 = $d.Paragraphs(47).Range.Text
$clickIdx = $full47.IndexOf("click")
$afterClickStart = $d.Paragraphs(47).Range.Start + $clickIdx + 5
$afterClickEnd = $d.Paragraphs(47).Range.End - 1
$afterClick = $d.Range($afterClickStart, $afterClickEnd)
$afterClick.Text = " en el botón de eliminar que está al lado del botón de edición."

# Paragraph 48 (third "Pasos" bullet)
# "Aparece un mensaje de confirmación y se le indica que la película ya
#  existe." -> "Aparece un mensaje de confirmación."
$p48 = $d.Paragraphs(48)
$r48 = $p48.Range
$r48.MoveEnd(1, -1) | Out-Null
$r48.Text = "Aparece un mensaje de confirmación."

# New bullet appended after paragraph 48, inheriting the same list
# style/numbering ("Prrafodelista", numId 8).
$p48b = $d.Paragraphs(48)
$p48b.Range.InsertParagraphAfter()
$p49 = $d.Paragraphs(49)
$r49 = $p49.Range
$r49.MoveEnd(1, -1) | Out-Null
$r49.Text = "Al afirmar la acción la película desaparece de la lista y de la base de datos."
